$wb = $excel.ActiveWorkbook

$newBase = "b095d723-d7d4-4dc8-96cf-d41fd8da1740"
$mdName = "$newBase.md"
$mdPath = "e2e\$newBase.md"
$zhXlf = "$newBase.6bfb1e39b49dcd2ab32c5c5e6f422e85ba37ccd2.zh-cn.xlf"
$deXlf = "$newBase.6bfb1e39b49dcd2ab32c5c5e6f422e85ba37ccd2.de-de.xlf"

$srcSha = "d169c5314dcbca5a5cfd38128123c0d98eccf3df"
$zhSha = "f48edd0bc64b7a7814122d3b69785b4257f7fa95"
$deSha = "47f7b6091567f86980cb2c49685214af664101b6"

$srcUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcSha/e2e/$mdName"
$zhUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$zhSha/e2e/$mdName"
$deUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$deSha/e2e/$mdName"

$dateFmt = "yyyy-mm-dd HH:mm:ss"
$hlColor = 6591981

# ---------------------------------------------------------------
# Sheet "Overview" (sheet1 / table3) - add row 4
# ---------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A4").Value = $mdName
$wsOv.Range("C4").Value = ".md"
$wsOv.Range("E4").Value = "Handed back: in sync with en-US"
$wsOv.Range("F4").Value = "Handed back: in sync with en-US"

$wsOv.Range("G4").NumberFormat = $dateFmt
$wsOv.Range("G4").Value = "2016-09-03 18:57:20"

$wsOv.Hyperlinks.Add($wsOv.Range("B4"), $srcUrl, "", "", $mdPath) | Out-Null
$wsOv.Range("B4").Font.Color = $hlColor

# ---------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table1) - add row 4
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("J4").Value = $zhXlf

$wsZh.Range("H4").NumberFormat = $dateFmt
$wsZh.Range("H4").Value = "2016-09-03 18:57:15"
$wsZh.Range("K4").NumberFormat = $dateFmt
$wsZh.Range("K4").Value = "2016-09-03 18:57:40"

$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $srcUrl, "", "", $mdName) | Out-Null
$wsZh.Range("A4").Font.Color = $hlColor
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $zhUrl, "", "", $mdName) | Out-Null
$wsZh.Range("I4").Font.Color = $hlColor

# ---------------------------------------------------------------
# Sheet "de-de" (sheet3 / table2) - add row 4
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("J4").Value = $deXlf

$wsDe.Range("H4").NumberFormat = $dateFmt
$wsDe.Range("H4").Value = "2016-09-03 18:57:20"
$wsDe.Range("K4").NumberFormat = $dateFmt
$wsDe.Range("K4").Value = "2016-09-03 18:57:47"

$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $srcUrl, "", "", $mdName) | Out-Null
$wsDe.Range("A4").Font.Color = $hlColor
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $deUrl, "", "", $mdName) | Out-Null
$wsDe.Range("I4").Font.Color = $hlColor
